$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.561.18'
$ws.Range("E2").Value = '  +2.59%  '
$ws.Range("D3").Value = '1.876.01'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  +1.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.37'
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5099'
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3924'
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08397'
$ws.Range("E9").Value = '  +3.57%  '
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.77'
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.268'
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("D13").Value = '1.873.58'
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("E14").Value = '  +2.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.268'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001107'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.42'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06736'
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.73'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.961'
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("D23").Value = '28.600.95'
$ws.Range("E23").Value = '  +2.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.246'
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").Value = '2.096.40'
$ws.Range("E26").Value = '  +3.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.91'
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.76'
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.367'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.58'
$ws.Range("E30").Value = '  +1.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.055'
$ws.Range("E32").Value = '  +2.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.802'
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.610'
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02457'
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06537'
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2186'
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.895'
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.263'
$ws.Range("E39").Value = '  +3.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.195'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6464'
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("E42").Value = '  +3.63%  '
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.008'
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6066'
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.08'
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.700'
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.029'
$ws.Range("E48").Value = '  +3.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.219'
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.33'
$ws.Range("E50").Value = '  +2.04%  '
$ws.Range("E51").Value = '  -5.48%  '
